$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44467
$ws.Range("J2").Value = 40
$ws.Range("O2").Value = "Región del Maule"
$ws.Range("D3").Value = 44369
$ws.Range("J3").Value = 60
$ws.Range("O3").Value = "Región Metropolitana"
$ws.Range("D5").Value = 44389
$ws.Range("J5").Value = 55
$ws.Range("O5").Value = "Región del Maule"
$ws.Range("D6").Value = 44348
$ws.Range("J6").Value = 35
$ws.Range("K6").Value = 7000
$ws.Range("L6").Value = 7000
$ws.Range("M6").Value = 7000
$ws.Range("O6").Value = "Región del Maule"
$ws.Range("P6").Value = 438
$ws.Range("D7").Value = 44312
$ws.Range("J7").Value = 40
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 438
$ws.Range("D8").Value = 44398
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 438
$ws.Range("D9").Value = 44420
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 45
$ws.Range("K9").Value = 8000
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 8000
$ws.Range("O9").Value = "Región Metropolitana"
$ws.Range("P9").Value = 500
$ws.Range("D10").Value = 44362
$ws.Range("J10").Value = 25
$ws.Range("K10").Value = 8000
$ws.Range("L10").Value = 8000
$ws.Range("M10").Value = 8000
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 500
$ws.Range("D11").Value = 44305
$ws.Range("J11").Value = 35
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("O11").Value = "Región del Maule"
$ws.Range("P11").Value = 438
$ws.Range("D12").Value = 44308
$ws.Range("J12").Value = 75
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 5000
$ws.Range("M12").Value = 5000
$ws.Range("O12").Value = "Región del Maule"
$ws.Range("P12").Value = 312
$ws.Range("D13").Value = 44354
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 8000
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = 8500
$ws.Range("O13").Value = "Región Metropolitana"
$ws.Range("P13").Value = 531
$ws.Range("D14").Value = 44354
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 9000
$ws.Range("L14").Value = 9000
$ws.Range("M14").Value = 9000
$ws.Range("P14").Value = 562
$ws.Range("D15").Value = 44397
$ws.Range("J15").Value = 40
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 500
$ws.Range("D16").Value = 44371
$ws.Range("J16").Value = 200
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("P16").Value = 438
$ws.Range("D17").Value = 44355
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 8000
$ws.Range("L17").Value = 8000
$ws.Range("M17").Value = 8000
$ws.Range("P17").Value = 500
$ws.Range("D18").Value = 44403
$ws.Range("J18").Value = 35
$ws.Range("K18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("M18").Value = 5000
$ws.Range("P18").Value = 312
$ws.Range("D19").Value = 44315
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 7000
$ws.Range("M19").Value = 7000
$ws.Range("P19").Value = 438
$ws.Range("D20").Value = 44386
$ws.Range("J20").Value = 40
$ws.Range("K20").Value = 7000
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = 7000
$ws.Range("O20").Value = "Región del Maule"
$ws.Range("P20").Value = 438
$ws.Range("D21").Value = 44313
$ws.Range("J21").Value = 20
$ws.Range("K21").Value = 7000
$ws.Range("L21").Value = 7000
$ws.Range("M21").Value = 7000
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 438
$ws.Range("D22").Value = 44396
$ws.Range("J22").Value = 80
$ws.Range("O22").Value = "Región Metropolitana"
$ws.Range("D23").Value = 44399
$ws.Range("J23").Value = 80
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("D24").Value = 44372
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 6000
$ws.Range("M24").Value = 6400
$ws.Range("P24").Value = 400
$ws.Range("D25").Value = 44392
$ws.Range("J25").Value = 95
$ws.Range("D26").Value = 44314
$ws.Range("I26").Value = "Segunda"
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 5000
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = 5000
$ws.Range("P26").Value = 312
